$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 103, shifting existing rows 103:167 down to 104:168.
$ws.Rows("103:103").Insert()

# Populate the new row 103 with the new daily price record.
$ws.Range("A103").Value = 8
$ws.Range("B103").Value = "Terminal La Palmera de La Serena"
$ws.Range("C103").Value = "Coquimbo"
$ws.Range("D103").Value = 44438
$ws.Range("D103").NumberFormat = $ws.Range("D104").NumberFormat
$ws.Range("E103").Value = 4
$ws.Range("F103").Value = 100112032
$ws.Range("G103").Value = "Zapallo italiano"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 600
$ws.Range("K103").Value = 13000
$ws.Range("L103").Value = 13500
$ws.Range("M103").Value = 13250
$ws.Range("N103").Value = "$/caja 50 unidades"
$ws.Range("O103").Value = "Región de Arica y Parinacota"
$ws.Range("P103").Value = 265
$ws.Range("Q103").Value = 50
$ws.Range("R103").Value = "Hortaliza"
